# Add a new "2021" column (O) to the dataset, mirroring the formatting of
# the existing "2020" column (N).
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: trailing border-only cell, same style as N2
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial($xlPasteFormats)

# Row 3: header year value, same style as N3
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial($xlPasteFormats)
$ws.Range("O3").Value = 2021

# Row 4: computed ratio (formula), same style as N4
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial($xlPasteFormats)
$ws.Range("O4").Formula = "=O5/O6*1000"

# Row 5: "removal of solid household waste" value (General number format,
# matching the sibling text cells on this row rather than the "0.00" used
# further left in the row)
$ws.Range("B5").Copy()
$ws.Range("O5").PasteSpecial($xlPasteFormats)
$ws.Range("O5").Value = 1229.5999999999999

# Row 6: "average annual population" value, same style as N6
$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial($xlPasteFormats)
$ws.Range("O6").Value = 6436.9

$excel.CutCopyMode = $false

# Move the active selection, as recorded after the edit
$ws.Range("P16").Select() | Out-Null
